# Sort the "Highscores" sheet data (A1:L14, header in row 1) by column K
# (cross-validation mean) in descending order, then leave the active
# selection on cell I6 - matching the manual "sort table by best score"
# pass that produced the printable poster figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Highscores")
$ws.Activate()

$sortRange = $ws.Range("A1:L14")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("K2:K14"), 0, 2, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

$ws.Range("I6").Select()
